$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [char]0x2083

$ws.Range("D2").Value = "'68.502.45"
$ws.Range("E2").Value = '  -6.58%  '

$ws.Range("D3").Value = "'3.742.87"
$ws.Range("E3").Value = '  -5.85%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = "'583.39"
$ws.Range("E5").Value = '  -5.61%  '

$ws.Range("D6").Value = "'180.59"
$ws.Range("E6").Value = '  +7.01%  '

$ws.Range("D7").Value = "'3.732.23"
$ws.Range("E7").Value = '  -5.85%  '

$ws.Range("D8").Value = "'0.638"
$ws.Range("E8").Value = '  -6.31%  '

$ws.Range("D9").Value = "'0.997"
$ws.Range("E9").Value = '  -0.24%  '

$ws.Range("D10").Value = "'0.725"
$ws.Range("E10").Value = '  -4.38%  '

$ws.Range("E11").Value = '  -10.19%  '

$ws.Range("D12").Value = "'54.19"
$ws.Range("E12").Value = '  -3.00%  '

$ws.Range("D13").Value = "'0.0000303"
$ws.Range("E13").Value = '  -9.91%  '

$ws.Range("D14").Value = "'10.88"
$ws.Range("E14").Value = '  -2.55%  '

$ws.Range("D15").Value = "'4.328.29"
$ws.Range("E15").Value = '  -6.17%  '

$ws.Range("D16").Value = "'3.730.64"
$ws.Range("E16").Value = '  -6.36%  '

$ws.Range("D17").Value = "'19.54"
$ws.Range("E17").Value = '  -4.32%  '

$ws.Range("D18").Value = "'13.18"
$ws.Range("E18").Value = '  -6.21%  '

$ws.Range("E19").Value = '  -6.94%  '

$ws.Range("E20").Value = '  -2.79%  '

$ws.Range("D21").Value = "'68.322.92"
$ws.Range("E21").Value = '  -6.56%  '

$ws.Range("D22").Value = "'413.47"
$ws.Range("E22").Value = '  -5.98%  '

$ws.Range("D23").Value = "'4.63"
$ws.Range("E23").Value = '  -5.09%  '

$ws.Range("D24").Value = "'89.26"
$ws.Range("E24").Value = '  -6.99%  '

$ws.Range("D25").Value = "'3.12"
$ws.Range("E25").Value = '  -7.36%  '

$ws.Range("E26").Value = '  -8.75%  '

$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = '  -0.75%  '

$ws.Range("D28").Value = "'3.87"
$ws.Range("E28").Value = '  -5.04%  '

$ws.Range("D29").Value = "'5.98"
$ws.Range("E29").Value = '  +0.23%  '

$ws.Range("D30").Value = "'9.69"
$ws.Range("E30").Value = '  -7.98%  '

$ws.Range("D31").Value = "'8.17"
$ws.Range("E31").Value = '  +4.49%  '

$ws.Range("D32").Value = "'33.30"
$ws.Range("E32").Value = '  -7.85%  '

$ws.Range("D33").Value = "'12.84"
$ws.Range("E33").Value = '  -6.17%  '

$ws.Range("E34").Value = '  -7.87%  '

$ws.Range("D35").Value = "'44.80"
$ws.Range("E35").Value = '  -6.40%  '

$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = "'616.63"
$ws.Range("E36").Value = '  -4.21%  '

$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = "'66.12"
$ws.Range("E37").Value = '  -6.89%  '

$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = ('0.0{0}0936' -f $sub3)
$ws.Range("E38").Value = '  -10.57%  '

$ws.Range("D39").Value = "'0.406"
$ws.Range("E39").Value = '  -5.60%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").Value = "'3.35"
$ws.Range("E40").Value = '  +8.14%  '

$ws.Range("B41").Value = 'Dai'
$ws.Range("C41").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = '  -0.21%  '

$ws.Range("D43").Value = "'0.138"
$ws.Range("E43").Value = '  -5.15%  '

$ws.Range("D44").Value = "'3.12"
$ws.Range("E44").Value = '  -8.06%  '

$ws.Range("E45").Value = '  -7.71%  '

$ws.Range("E46").Value = '  +3.50%  '

$ws.Range("D47").Value = "'9.65"
$ws.Range("E47").Value = '  -8.64%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = "'2.76"
$ws.Range("E48").Value = '  -14.36%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = "'0.136"
$ws.Range("E49").Value = '  -8.19%  '

$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = "'3.19"
$ws.Range("E50").Value = '  -6.45%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = "'2.756.68"
$ws.Range("E51").Value = '  -2.32%  '
